$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell F48 = "DONE"
$ws.Range("F48").Value = "DONE"

# New row 51
$ws.Range("A51").Value = "7/8/2025(Remote)"
$ws.Range("B51").Value = "Car Tracking Project"
$ws.Range("C51").Value = "Remember to change the extraction and limit rate to the normal"

# New row 52
$ws.Range("A52").Value = "7/8/2025(Remote)"
$ws.Range("B52").Value = "Car Tracking Project"
$ws.Range("C52").Value = "Make LinkedIn Post"

# New row 53
$ws.Range("C53").Value = "دراسة جدوى"

# New row 54
$ws.Range("C54").Value = "secret env variab;es"

# Update the sheet view: scroll position + selection
$ws.Range("B49").Select()
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 2
